# Auto-generated edit script applying cached-value updates to the
# Leve profit tables (currentAveragePrice / LevePrice / LeveProfit columns)
# across all 8 sheets, per the commit's scheduled-runner data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 4247.143
$ws.Range("J40").Value = 5404.8335
$ws.Range("L40").Value = 5404.8335
$ws.Range("N40").Value = -5754.8335
$ws.Range("H51").Value = 15874.821
$ws.Range("I51").Value = 6227.091
$ws.Range("J51").Value = 22117.47
$ws.Range("K51").Value = 6227.091
$ws.Range("L51").Value = 22117.47
$ws.Range("M51").Value = -5743.091
$ws.Range("N51").Value = -23085.47
$ws.Range("H103").Value = 1326.8182
$ws.Range("I103").Value = 650
$ws.Range("J103").Value = 1477.2222
$ws.Range("K103").Value = 1950
$ws.Range("L103").Value = 4431.6666
$ws.Range("M103").Value = -1364
$ws.Range("N103").Value = -5603.6666
$ws.Range("H129").Value = 2120.4
$ws.Range("I129").Value = 1526.6666
$ws.Range("J129").Value = 2606.182
$ws.Range("K129").Value = 4579.9998
$ws.Range("L129").Value = 7818.545999999999
$ws.Range("M129").Value = 420.0002000000004
$ws.Range("N129").Value = -17818.546
$ws.Range("H132").Value = 2014.75
$ws.Range("I132").Value = 1187
$ws.Range("K132").Value = 3561
$ws.Range("M132").Value = -1031
$ws.Range("H135").Value = 2345.5652
$ws.Range("I135").Value = 1229.2222
$ws.Range("K135").Value = 11062.9998
$ws.Range("M135").Value = -8527.9998
$ws.Range("H136").Value = 119995
$ws.Range("J136").Value = 119995
$ws.Range("L136").Value = 119995
$ws.Range("N136").Value = -130195
$ws.Range("H138").Value = 2917.7778
$ws.Range("I138").Value = 3360.8572
$ws.Range("K138").Value = 10082.5716
$ws.Range("M138").Value = -4942.571599999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 33320.92
$ws.Range("I32").Value = 35129.086
$ws.Range("K32").Value = 35129.086
$ws.Range("M32").Value = -34842.086
$ws.Range("H61").Value = 13031.588
$ws.Range("I61").Value = 12047.546
$ws.Range("K61").Value = 12047.546
$ws.Range("M61").Value = -11835.546
$ws.Range("H136").Value = 13031.588
$ws.Range("I136").Value = 12047.546
$ws.Range("K136").Value = 36142.638
$ws.Range("M136").Value = -33592.638
$ws.Range("H141").Value = 120000
$ws.Range("J141").Value = 120000
$ws.Range("L141").Value = 120000
$ws.Range("N141").Value = -130360

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2841.64
$ws.Range("I20").Value = 2502.1177
$ws.Range("J20").Value = 3563.125
$ws.Range("K20").Value = 2502.1177
$ws.Range("L20").Value = 3563.125
$ws.Range("M20").Value = -2255.1177
$ws.Range("N20").Value = -4057.125
$ws.Range("H107").Value = 1312.2727
$ws.Range("I107").Value = 1260.0714
$ws.Range("J107").Value = 1403.625
$ws.Range("K107").Value = 1260.0714
$ws.Range("L107").Value = 1403.625
$ws.Range("M107").Value = 659.9286
$ws.Range("N107").Value = -5243.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1312.2222
$ws.Range("I105").Value = 830.1429000000001
$ws.Range("K105").Value = 830.1429000000001
$ws.Range("M105").Value = 916.8570999999999
$ws.Range("H108").Value = 120000
$ws.Range("J108").Value = 120000
$ws.Range("L108").Value = 120000
$ws.Range("N108").Value = -127680
$ws.Range("H132").Value = 19255.395
$ws.Range("I132").Value = 1803.3572
$ws.Range("K132").Value = 5410.071599999999
$ws.Range("M132").Value = -2880.071599999999
$ws.Range("H134").Value = 3034.4146
$ws.Range("I134").Value = 2403.3125
$ws.Range("J134").Value = 5278.3335
$ws.Range("K134").Value = 7209.9375
$ws.Range("L134").Value = 15835.0005
$ws.Range("M134").Value = -4674.9375
$ws.Range("N134").Value = -20905.0005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 126940.81
$ws.Range("I68").Value = 401200.6
$ws.Range("K68").Value = 1203601.8
$ws.Range("M68").Value = -1202790.8
$ws.Range("H71").Value = 126940.81
$ws.Range("I71").Value = 401200.6
$ws.Range("K71").Value = 3610805.4
$ws.Range("M71").Value = -3606749.4
$ws.Range("H113").Value = 1850.75
$ws.Range("J113").Value = 2334.3333
$ws.Range("L113").Value = 7002.999899999999
$ws.Range("N113").Value = -11342.9999
$ws.Range("H131").Value = 8549839
$ws.Range("J131").Value = 5391.647
$ws.Range("L131").Value = 16174.941
$ws.Range("N131").Value = -26254.941

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7398.375
$ws.Range("I70").Value = 8247
$ws.Range("K70").Value = 8247
$ws.Range("M70").Value = -7977
$ws.Range("H73").Value = 7398.375
$ws.Range("I73").Value = 8247
$ws.Range("K73").Value = 8247
$ws.Range("M73").Value = -7311
$ws.Range("H102").Value = 1346.4375
$ws.Range("I102").Value = 1369.5333
$ws.Range("J102").Value = 1000
$ws.Range("K102").Value = 1369.5333
$ws.Range("L102").Value = 1000
$ws.Range("M102").Value = 252.4666999999999
$ws.Range("N102").Value = -4244
$ws.Range("H107").Value = 376.82608
$ws.Range("I107").Value = 375.35715
$ws.Range("K107").Value = 375.35715
$ws.Range("M107").Value = 1544.64285
$ws.Range("H122").Value = 10740.714
$ws.Range("I122").Value = 800
$ws.Range("J122").Value = 12397.5
$ws.Range("K122").Value = 2400
$ws.Range("L122").Value = 37192.5
$ws.Range("M122").Value = 50
$ws.Range("N122").Value = -42092.5
$ws.Range("H126").Value = 3388.3333
$ws.Range("I126").Value = 2480.3076
$ws.Range("J126").Value = 4863.875
$ws.Range("K126").Value = 7440.9228
$ws.Range("L126").Value = 14591.625
$ws.Range("M126").Value = -4970.9228
$ws.Range("N126").Value = -19531.625
$ws.Range("H132").Value = 5110.2554
$ws.Range("I132").Value = 3790.0715
$ws.Range("K132").Value = 11370.2145
$ws.Range("M132").Value = -8840.2145

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5565.7437
$ws.Range("I40").Value = 5219.6763
$ws.Range("J40").Value = 7919
$ws.Range("K40").Value = 5219.6763
$ws.Range("L40").Value = 7919
$ws.Range("M40").Value = -5083.6763
$ws.Range("N40").Value = -8191
$ws.Range("H61").Value = 2575.9375
$ws.Range("I61").Value = 1436.909
$ws.Range("K61").Value = 1436.909
$ws.Range("M61").Value = -1234.909
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H113").Value = 2575.9375
$ws.Range("I113").Value = 1436.909
$ws.Range("K113").Value = 1436.909
$ws.Range("M113").Value = 733.0909999999999
$ws.Range("H133").Value = 111111
$ws.Range("J133").Value = 111111
$ws.Range("L133").Value = 111111
$ws.Range("N133").Value = -116171

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H82").Value = 50301
$ws.Range("J82").Value = 50301
$ws.Range("L82").Value = 50301
$ws.Range("N82").Value = -51067
$ws.Range("H85").Value = 50301
$ws.Range("J85").Value = 50301
$ws.Range("L85").Value = 50301
$ws.Range("N85").Value = -52953
$ws.Range("H86").Value = 78000
$ws.Range("J86").Value = 78000
$ws.Range("L86").Value = 78000
$ws.Range("N86").Value = -80246
$ws.Range("H89").Value = 78000
$ws.Range("J89").Value = 78000
$ws.Range("L89").Value = 390000
$ws.Range("N89").Value = -401232
$ws.Range("H110").Value = 51644
$ws.Range("J110").Value = 51644
$ws.Range("L110").Value = 51644
$ws.Range("N110").Value = -59824
$ws.Range("H111").Value = 69999
$ws.Range("J111").Value = 69999
$ws.Range("L111").Value = 69999
$ws.Range("N111").Value = -78179
$ws.Range("H132").Value = 3716.6667
$ws.Range("I132").Value = 2394.8386
$ws.Range("K132").Value = 7184.5158
$ws.Range("M132").Value = -4654.5158
$ws.Range("H136").Value = 3573.6428
$ws.Range("I136").Value = 1583.5
$ws.Range("K136").Value = 4750.5
$ws.Range("M136").Value = -2200.5

